# Commit: "base model changed to efficientnet b4, batch size reduced to 16"
#
# The experiment log lives on the "Тесты" sheet (sheet1 / the active sheet).
# Row 21 (test #20) gets its "Результаты" (results) cell filled in, and two
# brand-new test rows are appended: row 22 (test #21) and row 23 (test #22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21 (test #20): record the IoU results for that run --------------
$ws.Range("H21").Value = "Train IoU: 0.46, Val IoU: 0.46. "

# --- Row 22 (test #21): new base-model / loss-function experiment --------
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 40
$ws.Range("D22").Value = 14
$ws.Range("F22").Value = "Функция потерь заменена на комбинацию MSE и расстояния между центрами (50/50)"
$ws.Range("G22").Value = "параметры теста 4"
$ws.Range("H22").Value = "Train IoU: 0.54, Val IoU: 0.54. "

# Row grew taller to fit the wrapped "changes" text.
$ws.Rows.Item(22).RowHeight = 45

# --- Row 23 (test #22): follow-up run, proportion tweak -------------------
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 40
$ws.Range("D23").Value = 14
$ws.Range("F23").Value = "MSE/distance = 30/70"
$ws.Range("G23").Value = "параметры теста 4"

# Commit id is filled in last (after the rows were drafted), landing after
# the row-23 strings in the shared-string table.
$ws.Range("I22").Value = "4df1807"

# --- View state: scroll down a bit and leave the selection on J23 --------
$aw = $excel.ActiveWindow
$aw.ScrollRow = 13
$aw.ScrollColumn = 2
$ws.Range("J23").Select()
